# Add/update metadata report for Akurana
# Appends a new row (row 5) to the "Metadata Report" sheet, mirroring the
# existing DEC-2024 placeholder row (row 4) but extended out to column AY.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata Report")

# Leading descriptive columns (A-E)
$ws.Cells.Item(5, 1).Value = 2024
$ws.Cells.Item(5, 2).Value = "DEC"
$ws.Cells.Item(5, 3).Value = "31/12-01/12"
$ws.Cells.Item(5, 4).Value = "-"
$ws.Cells.Item(5, 5).Value = "Akurana"

# Remaining columns (F through AY) are all placeholder "-" values.
for ($col = 6; $col -le 51; $col++) {
    $ws.Cells.Item(5, $col).Value = "-"
}
